# Apply the edit described by the diff to the ModuleController worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModuleController")

# Update the shared-string values that changed:
#   B8: FASB2             -> sample
#   B2: PortFolioInsight -> Modul1
#   B3: FASB              -> Module2
# (order matters for how Excel rebuilds the shared-string table on save)
$ws.Range("B8").Value = "sample"
$ws.Range("B2").Value = "Modul1"
$ws.Range("B3").Value = "Module2"

# Update the selected cell/range in the sheet view (was A4, now B9)
$ws.Range("B9").Select()
